# Energy Consumption per Filter.xlsx
# "chaged 1900 to MERV 13"
#
# The filters previously labeled "1900-1", "1900-2", "1900-4", "1900-5" are
# renamed: the old "FG-*" filters become "MERV 13-*", and the old "1900-*"
# filters become the new "FG-*". The underlying measured values stay
# attached to the (now renamed) filter identity, which means - because the
# rows keep their original order/position in the table - the energy
# readings for rows 2/4 and 3/5 (and, in the second block, 12/14 and 13/15)
# swap places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Filter ID" labels (column C) -----------------------------
$ws.Range("C2").Value = "FG-1"
$ws.Range("C3").Value = "FG-2"
$ws.Range("C4").Value = "MERV 13-1"
$ws.Range("C5").Value = "MERV 13-2"
$ws.Range("C8").Value = "MERV 13-4"
$ws.Range("C9").Value = "MERV 13-5"

$ws.Range("C12").Value = "FG-1"
$ws.Range("C13").Value = "FG-2"
$ws.Range("C14").Value = "MERV 13-1"
$ws.Range("C15").Value = "MERV 13-2"
$ws.Range("C18").Value = "MERV 13-4"
$ws.Range("C19").Value = "MERV 13-5"

# --- Swap the energy readings (columns D:G) that travel with the renamed
#     filters (rows 2<->4, 3<->5, 12<->14, 13<->15) -------------------------
$ws.Range("D2").Value = 6.2184908458334
$ws.Range("E2").Value = 2.284984018611098
$ws.Range("F2").Value = 0.212494666666656
$ws.Range("G2").Value = 2.497478685277773

$ws.Range("D3").Value = 7.127588941944426
$ws.Range("E3").Value = 2.245040258333344
$ws.Range("F3").Value = 0.208120999999981
$ws.Range("G3").Value = 2.453161258333293

$ws.Range("D4").Value = 7.531882454722191
$ws.Range("E4").Value = 2.423465211111133
$ws.Range("F4").Value = 0.2240219999999769
$ws.Range("G4").Value = 2.647487211111092

$ws.Range("D5").Value = 7.390565565277755
$ws.Range("E5").Value = 2.353605584166654
$ws.Range("F5").Value = 0.2190806666666559
$ws.Range("G5").Value = 2.572686250833364

$ws.Range("D12").Value = 10.7667291222218
$ws.Range("E12").Value = 2.980704309444402
$ws.Range("F12").Value = 0.2444383333333209
$ws.Range("G12").Value = 3.225142642777744

$ws.Range("D13").Value = 10.43482538805555
$ws.Range("E13").Value = 2.843067705833336
$ws.Range("F13").Value = 0.2363996666666406
$ws.Range("G13").Value = 3.079467372500054

$ws.Range("D14").Value = 10.59388187749992
$ws.Range("E14").Value = 2.90196882166668
$ws.Range("F14").Value = 0.2475556666666464
$ws.Range("G14").Value = 3.14952448833322

$ws.Range("D15").Value = 9.927605422222198
$ws.Range("E15").Value = 2.603467600833405
$ws.Range("F15").Value = 0.2272253333333117
$ws.Range("G15").Value = 2.830692934166631
